# Horarios actualizados Línea 141 - 135
# Applies the 07:38:39 scrape update to all three sheets of the workbook.

$wb = $excel.ActiveWorkbook

$newTimestamp = "07:38:39"

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTimestamp"
$ws1.Range("A3").Value = "Total filas: 48"

# Rows 40/41 had their "Linea" values swapped in this update.
$ws1.Cells.Item(40, 3).Value = "15_ABASTO"
$ws1.Cells.Item(41, 3).Value = "11_ETCHEVERRY"

# New scrape rows appended at the bottom (48-53).
$sheet1NewRows = @(
    @("07:38:39", "08:58", "215A_EL PATO", 80,  "LP1912"),
    @("07:38:39", "09:06", "16_SANTA ANA", 88,  "LP1912"),
    @("07:38:39", "09:16", "27_EL RETIRO", 98,  "LP1912"),
    @("07:38:39", "09:17", "14_ABASTO",    99,  "LP1912"),
    @("07:38:39", "09:18", "15X38_ABASTO", 100, "LP1912"),
    @("07:38:39", "09:29", "10_OLMOS",     111, "LP1912")
)

$r = 48
foreach ($row in $sheet1NewRows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTimestamp"
$ws2.Range("A3").Value = "Total filas: 13"

$ws2.Cells.Item(18, 1).Value = "07:38:39"
$ws2.Cells.Item(18, 2).Value = "08:58"
$ws2.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(18, 4).Value = 80
$ws2.Cells.Item(18, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTimestamp"
